$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing "DuréeVie" values (column C) into the new column E,
# then update column C with the recalculated values from the new algorithm.
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 3).Value2
}

$ws.Range("C2").Value = 20
$ws.Range("C6").Value = 20
$ws.Range("C9").Value = 4
$ws.Range("C11").Value = 5
$ws.Range("C16").Value = 20

$ws.Range("H8").Select() | Out-Null
